$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "27.921.86"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.629.75"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.58%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "211.83"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  +0.03%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "23.34"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.20%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.256"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("E10").Value = "  -0.31%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0879"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.45%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.861.04"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.56%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.627.38"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("E15").Value = "  -2.47%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "65.58"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.63%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "27.920.35"
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "230.18"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.73%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.0₃0722"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.28%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.63"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.38%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("E23").Value = "  -5.10%  "
$ws.Range("E24").Value = "  -1.19%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "154.76"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.99%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "6.93"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("E28").Value = "  -1.25%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("E32").Value = "  +2.16%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.08"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.78%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.401.50"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("E36").Value = "  +11.11%  "
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("E38").Value = "  +1.94%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.555"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("E40").Value = "  -3.32%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.02"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("E46").Value = "  -0.46%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.771.34"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.47%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "87.97"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("E49").Value = "  -2.86%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("E51").Value = "  -0.16%  "
